$d = $word.ActiveDocument
$d.Content.Find.Execute("816÷9=90, 6", $true, $false, $false, $false, $false, $true, 1, $false, "163÷7=23, 2", 2) | Out-Null
$d.Content.Find.Execute("382÷3=127, 1", $true, $false, $false, $false, $false, $true, 1, $false, "910÷3=303, 1", 2) | Out-Null
$d.Content.Find.Execute("144÷7=20, 4", $true, $false, $false, $false, $false, $true, 1, $false, "915÷3=305, 0", 2) | Out-Null
$d.Content.Find.Execute("562÷6=93, 4", $true, $false, $false, $false, $false, $true, 1, $false, "434÷2=217, 0", 2) | Out-Null
$d.Content.Find.Execute("156÷8=19, 4", $true, $false, $false, $false, $false, $true, 1, $false, "755÷5=151, 0", 2) | Out-Null
$d.Content.Find.Execute("481÷4=120, 1", $true, $false, $false, $false, $false, $true, 1, $false, "185÷9=20, 5", 2) | Out-Null
$d.Content.Find.Execute("403÷5=80, 3", $true, $false, $false, $false, $false, $true, 1, $false, "341÷9=37, 8", 2) | Out-Null
$d.Content.Find.Execute("766÷4=191, 2", $true, $false, $false, $false, $false, $true, 1, $false, "563÷6=93, 5", 2) | Out-Null
$d.Content.Find.Execute("600÷3=200, 0", $true, $false, $false, $false, $false, $true, 1, $false, "158÷4=39, 2", 2) | Out-Null
$d.Content.Find.Execute("560÷9=62, 2", $true, $false, $false, $false, $false, $true, 1, $false, "392÷3=130, 2", 2) | Out-Null
$d.Content.Find.Execute("715÷9=79, 4", $true, $false, $false, $false, $false, $true, 1, $false, "552÷2=276, 0", 2) | Out-Null
$d.Content.Find.Execute("438÷5=87, 3", $true, $false, $false, $false, $false, $true, 1, $false, "705÷7=100, 5", 2) | Out-Null
$d.Content.Find.Execute("935÷7=133, 4", $true, $false, $false, $false, $false, $true, 1, $false, "431÷4=107, 3", 2) | Out-Null
$d.Content.Find.Execute("901÷8=112, 5", $true, $false, $false, $false, $false, $true, 1, $false, "523÷6=87, 1", 2) | Out-Null
$d.Content.Find.Execute("417÷3=139, 0", $true, $false, $false, $false, $false, $true, 1, $false, "963÷3=321, 0", 2) | Out-Null
$d.Content.Find.Execute("613÷4=153, 1", $true, $false, $false, $false, $false, $true, 1, $false, "174÷5=34, 4", 2) | Out-Null
$d.Content.Find.Execute("711÷6=118, 3", $true, $false, $false, $false, $false, $true, 1, $false, "482÷9=53, 5", 2) | Out-Null
$d.Content.Find.Execute("982÷7=140, 2", $true, $false, $false, $false, $false, $true, 1, $false, "753÷4=188, 1", 2) | Out-Null
$d.Content.Find.Execute("673÷2=336, 1", $true, $false, $false, $false, $false, $true, 1, $false, "763÷4=190, 3", 2) | Out-Null
$d.Content.Find.Execute("618÷8=77, 2", $true, $false, $false, $false, $false, $true, 1, $false, "690÷3=230, 0", 2) | Out-Null
$d.Content.Find.Execute("749÷9=83, 2", $true, $false, $false, $false, $false, $true, 1, $false, "297÷2=148, 1", 2) | Out-Null
$d.Content.Find.Execute("438÷2=219, 0", $true, $false, $false, $false, $false, $true, 1, $false, "952÷7=136, 0", 2) | Out-Null
$d.Content.Find.Execute("471÷4=117, 3", $true, $false, $false, $false, $false, $true, 1, $false, "774÷6=129, 0", 2) | Out-Null
$d.Content.Find.Execute("366÷3=122, 0", $true, $false, $false, $false, $false, $true, 1, $false, "320÷7=45, 5", 2) | Out-Null
$d.Content.Find.Execute("674÷7=96, 2", $true, $false, $false, $false, $false, $true, 1, $false, "710÷9=78, 8", 2) | Out-Null
